# GUI: Updated the statistics.
#
# The "Running" row (row 5) of Sheet1 pulls its numbers from an external
# workbook link ([4]Sheet1 -> Running/_Test_Suite_Statistics.xlsx):
#   C5 = [4]Sheet1!$G$2   (Automated Test Suites)
#   D5 = [4]Sheet1!$E$1   (Total Test Suites)
#   G5 = [4]Sheet1!$G$5   (Automated Test Cases)
#   H5 = [4]Sheet1!$G$4   (Total Test Cases)
# The external source's numbers were updated (as if the linked workbook
# had been edited and the link refreshed), so those four cached values
# change here. Every other statistic on the sheet (L1, N1, N2, L3, N3,
# P3, E5, I5, L5, L6, L7, ...) is a plain in-workbook formula, so simply
# recalculates on its own once the source numbers above change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value2 = 5
$ws.Range("D5").Value2 = 11
$ws.Range("G5").Value2 = 11
$ws.Range("H5").Value2 = 56
